$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve existing formatting for column D, force text entry so numeric-looking
# strings (e.g. "1.007") are stored as text, matching the original inline strings,
# then restore the original style so no stray formatting is introduced.
$dRange = $ws.Range("D2:D51")
$origDStyle = $dRange.Style
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '27.366.93'
$ws.Range("D3").Value = '1.715.56'
$ws.Range("D4").Value = '1.007'
$ws.Range("D5").Value = '224.93'
$ws.Range("D6").Value = '0.5280'
$ws.Range("D8").Value = '0.06671'
$ws.Range("D10").Value = '20.81'
$ws.Range("D11").Value = '0.07763'
$ws.Range("D12").Value = '4.479'
$ws.Range("D13").Value = '1.952.74'
$ws.Range("D14").Value = '1.718.03'
$ws.Range("D15").Value = '0.5801'
$ws.Range("D16").Value = '0.0₅8201'
$ws.Range("D17").Value = '67.83'
$ws.Range("D18").Value = '27.383.64'
$ws.Range("D19").Value = '219.52'
$ws.Range("D20").Value = '1.009'
$ws.Range("D21").Value = '4.655'
$ws.Range("D22").Value = '10.44'
$ws.Range("D23").Value = '6.058'
$ws.Range("D25").Value = '145.24'
$ws.Range("D26").Value = '1.723'
$ws.Range("D27").Value = '0.1207'
$ws.Range("D28").Value = '7.228'
$ws.Range("D29").Value = '16.22'
$ws.Range("D30").Value = '0.05334'
$ws.Range("D31").Value = '1.296'
$ws.Range("D32").Value = '3.488'
$ws.Range("D33").Value = '3.397'
$ws.Range("D35").Value = '2.825'
$ws.Range("D36").Value = '0.9547'
$ws.Range("D37").Value = '2.401'
$ws.Range("D38").Value = '0.5893'
$ws.Range("D39").Value = '1.191.14'
$ws.Range("D40").Value = '0.01654'
$ws.Range("D41").Value = '5.823'
$ws.Range("D43").Value = '0.8416'
$ws.Range("D44").Value = '101.24'
$ws.Range("D45").Value = '1.859.59'
$ws.Range("D46").Value = '0.0₈118'
$ws.Range("D47").Value = '57.59'
$ws.Range("D48").Value = '0.4552'
$ws.Range("D49").Value = '1.007'
$ws.Range("D50").Value = '8.163'
$ws.Range("D51").Value = '0.05236'

$dRange.Style = $origDStyle

# Column E values are percentage strings (with padding spaces) and are stored as
# text natively since they are not valid numbers.
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("E3").Value = '  +0.30%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +2.24%  '
$ws.Range("E9").Value = '  +0.78%  '
$ws.Range("E10").Value = '  -0.85%  '
$ws.Range("E11").Value = '  +1.61%  '
$ws.Range("E12").Value = '  -1.72%  '
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("E16").Value = '  +0.62%  '
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("E18").Value = '  +0.39%  '
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("E20").Value = '  +0.30%  '
$ws.Range("E21").Value = '  -0.46%  '
$ws.Range("E22").Value = '  -0.35%  '
$ws.Range("E23").Value = '  +1.98%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  +1.65%  '
$ws.Range("E26").Value = '  -0.53%  '
$ws.Range("E27").Value = '  -0.79%  '
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("E29").Value = '  -0.23%  '
$ws.Range("E30").Value = '  -0.89%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("E33").Value = '  -0.17%  '
$ws.Range("E34").Value = '  +0.35%  '
$ws.Range("E35").Value = '  -1.50%  '
$ws.Range("E36").Value = '  +0.77%  '
$ws.Range("E37").Value = '  -1.18%  '
$ws.Range("E38").Value = '  +0.85%  '
$ws.Range("E39").Value = '  +14.83%  '
$ws.Range("E40").Value = '  +1.59%  '
$ws.Range("E41").Value = '  -0.76%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  +0.30%  '
$ws.Range("E44").Value = '  +0.21%  '
$ws.Range("E45").Value = '  +0.10%  '
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("E47").Value = '  -0.55%  '
$ws.Range("E48").Value = '  +1.07%  '
$ws.Range("E49").Value = '  +0.58%  '
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("E51").Value = '  -0.20%  '
